$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 24,14
$arr[0,0] = 11.73582396971045
$arr[0,1] = 8.198143628860628
$arr[0,2] = 5.742490069362258
$arr[0,3] = 11.71729909709081
$arr[0,4] = 27.40552900906781
$arr[0,5] = 0
$arr[0,6] = 7.344005520526261
$arr[0,7] = 24.69008135183563
$arr[0,8] = 0
$arr[0,9] = 9.102014730398251
$arr[0,10] = 9.757681587473192
$arr[0,11] = 0
$arr[0,12] = 0
$arr[0,13] = 24.62061181298066
$arr[1,0] = 11.44897535467174
$arr[1,1] = 8.169600342180345
$arr[1,2] = 5.69305807224514
$arr[1,3] = 11.73032061351807
$arr[1,4] = 27.4605122152918
$arr[1,5] = 0
$arr[1,6] = 7.344005520526261
$arr[1,7] = 24.78850513011825
$arr[1,8] = 0
$arr[1,9] = 8.895974339975714
$arr[1,10] = 9.742494873551111
$arr[1,11] = 0
$arr[1,12] = 0
$arr[1,13] = 24.70489543553914
$arr[2,0] = 11.27087887610328
$arr[2,1] = 8.151959755600826
$arr[2,2] = 5.662036351539757
$arr[2,3] = 11.74045278583213
$arr[2,4] = 27.50119698236927
$arr[2,5] = 0
$arr[2,6] = 7.344005520526261
$arr[2,7] = 24.85356120080131
$arr[2,8] = 0
$arr[2,9] = 8.768145275271289
$arr[2,10] = 9.734843055465978
$arr[2,11] = 0
$arr[2,12] = 0
$arr[2,13] = 24.76173537560008
$arr[3,0] = 11.19791374239924
$arr[3,1] = 8.144744233902859
$arr[3,2] = 5.649231290319371
$arr[3,3] = 11.74511918157891
$arr[3,4] = 27.51951383276444
$arr[3,5] = 0
$arr[3,6] = 7.344005520526261
$arr[3,7] = 24.88123389307233
$arr[3,8] = 0
$arr[3,9] = 8.715794003208801
$arr[3,10] = 9.732148133322765
$arr[3,11] = 0
$arr[3,12] = 0
$arr[3,13] = 24.78617573377926
$arr[4,0] = 11.18577761330035
$arr[4,1] = 8.143544540583811
$arr[4,2] = 5.647095297006776
$arr[4,3] = 11.74592649430654
$arr[4,4] = 27.5226601485653
$arr[4,5] = 0
$arr[4,6] = 7.344005520526261
$arr[4,7] = 24.88589907042705
$arr[4,8] = 0
$arr[4,9] = 8.707087573488375
$arr[4,10] = 9.73172627735762
$arr[4,11] = 0
$arr[4,12] = 0
$arr[4,13] = 24.79031113731163
$arr[5,0] = 11.26989627568565
$arr[5,1] = 8.151862550526968
$arr[5,2] = 5.661864313086904
$arr[5,3] = 11.74051354230501
$arr[5,4] = 27.5014369803176
$arr[5,5] = 0
$arr[5,6] = 7.344005520526261
$arr[5,7] = 24.85392970130923
$arr[5,8] = 0
$arr[5,9] = 8.767440202934887
$arr[5,10] = 9.734804993839058
$arr[5,11] = 0
$arr[5,12] = 0
$arr[5,13] = 24.76205981669777
$arr[6,0] = 11.63738960603847
$arr[6,1] = 8.188326981283438
$arr[6,2] = 5.72558806095596
$arr[6,3] = 11.7213455068645
$arr[6,4] = 27.4230478174426
$arr[6,5] = 0
$arr[6,6] = 7.344005520526261
$arr[6,7] = 24.72305779157516
$arr[6,8] = 0
$arr[6,9] = 9.031288079031089
$arr[6,10] = 9.752099641275681
$arr[6,11] = 0
$arr[6,12] = 0
$arr[6,13] = 24.64861520939219
$arr[7,0] = 12.3380106980313
$arr[7,1] = 8.25883297269994
$arr[7,2] = 5.845013524147194
$arr[7,3] = 11.70070291672803
$arr[7,4] = 27.32442202873927
$arr[7,5] = 0
$arr[7,6] = 7.344005520526261
$arr[7,7] = 24.50313892265632
$arr[7,8] = 0
$arr[7,9] = 9.535221206957969
$arr[7,10] = 9.799162481231969
$arr[7,11] = 0
$arr[7,12] = 0
$arr[7,13] = 24.46662766829284
$arr[8,0] = 12.83505702715021
$arr[8,1] = 8.309906016622728
$arr[8,2] = 5.929095737029087
$arr[8,3] = 11.69584752845386
$arr[8,4] = 27.28572048647932
$arr[8,5] = 0
$arr[8,6] = 7.344005520526261
$arr[8,7] = 24.36399854537254
$arr[8,8] = 0
$arr[8,9] = 9.893472513480548
$arr[8,10] = 9.84157015799191
$arr[8,11] = 0
$arr[8,12] = 0
$arr[8,13] = 24.35772293067212
$arr[9,0] = 13.05629483247685
$arr[9,1] = 8.332957245999275
$arr[9,2] = 5.966493746375204
$arr[9,3] = 11.69587015974643
$arr[9,4] = 27.2754686214209
$arr[9,5] = 0
$arr[9,6] = 7.344005520526261
$arr[9,7] = 24.30558206826068
$arr[9,8] = 0
$arr[9,9] = 10.05312446611084
$arr[9,10] = 9.86251765950013
$arr[9,11] = 0
$arr[9,12] = 0
$arr[9,13] = 24.31359021063861
$arr[10,0] = 13.13929349360037
$arr[10,1] = 8.341657745441566
$arr[10,2] = 5.98052812393517
$arr[10,3] = 11.69619866660762
$arr[10,4] = 27.27264478405369
$arr[10,5] = 0
$arr[10,6] = 7.344005520526261
$arr[10,7] = 24.28416390671157
$arr[10,8] = 0
$arr[10,9] = 10.11304854139376
$arr[10,10] = 9.870683780162508
$arr[10,11] = 0
$arr[10,12] = 0
$arr[10,13] = 24.2976580381323
$arr[11,0] = 13.12145409337973
$arr[11,1] = 8.339785241600591
$arr[11,2] = 5.977511315825139
$arr[11,3] = 11.69611370416239
$arr[11,4] = 27.27320587227348
$arr[11,5] = 0
$arr[11,6] = 7.344005520526261
$arr[11,7] = 24.28874540940967
$arr[11,8] = 0
$arr[11,9] = 10.1001673489312
$arr[11,10] = 9.86891473615286
$arr[11,11] = 0
$arr[11,12] = 0
$arr[11,13] = 24.30105459524005
$arr[12,0] = 13.06313922928525
$arr[12,1] = 8.333673625762501
$arr[12,2] = 5.967650945332053
$arr[12,3] = 11.69589078013835
$arr[12,4] = 27.27521508973613
$arr[12,5] = 0
$arr[12,6] = 7.344005520526261
$arr[12,7] = 24.30380588892454
$arr[12,8] = 0
$arr[12,9] = 10.05806543410751
$arr[12,10] = 9.863184832366628
$arr[12,11] = 0
$arr[12,12] = 0
$arr[12,13] = 24.31226381810301
$arr[13,0] = 13.02731597378005
$arr[13,1] = 8.329926304354856
$arr[13,2] = 5.9615944411377
$arr[13,3] = 11.69579586738462
$arr[13,4] = 27.27658363228064
$arr[13,5] = 0
$arr[13,6] = 7.344005520526261
$arr[13,7] = 24.3131224425092
$arr[13,8] = 0
$arr[13,9] = 10.03220585237899
$arr[13,10] = 9.859705408305649
$arr[13,11] = 0
$arr[13,12] = 0
$arr[13,13] = 24.3192314238212
$arr[14,0] = 12.82049371165015
$arr[14,1] = 8.308395702084018
$arr[14,2] = 5.926634126053282
$arr[14,3] = 11.69589087462449
$arr[14,4] = 27.28653853487016
$arr[14,5] = 0
$arr[14,6] = 7.344005520526261
$arr[14,7] = 24.3679144893998
$arr[14,8] = 0
$arr[14,9] = 9.882967272160059
$arr[14,10] = 9.840234144951792
$arr[14,11] = 0
$arr[14,12] = 0
$arr[14,13] = 24.36071613118304
$arr[15,0] = 12.69230921763508
$arr[15,1] = 8.295139081521505
$arr[15,2] = 5.904965440039645
$arr[15,3] = 11.6965200578225
$arr[15,4] = 27.29452983418524
$arr[15,5] = 0
$arr[15,6] = 7.344005520526261
$arr[15,7] = 24.4027782100435
$arr[15,8] = 0
$arr[15,9] = 9.790523539726582
$arr[15,10] = 9.828710275913762
$arr[15,11] = 0
$arr[15,12] = 0
$arr[15,13] = 24.38755250920039
$arr[16,0] = 12.61812670288042
$arr[16,1] = 8.28749721729873
$arr[16,2] = 5.892422338246235
$arr[16,3] = 11.69709204121077
$arr[16,4] = 27.29981837048067
$arr[16,5] = 0
$arr[16,6] = 7.344005520526261
$arr[16,7] = 24.42329007452901
$arr[16,8] = 0
$arr[16,9] = 9.737043064677055
$arr[16,10] = 9.822238204709523
$arr[16,11] = 0
$arr[16,12] = 0
$arr[16,13] = 24.40349704181805
$arr[17,0] = 12.59293427916783
$arr[17,1] = 8.284906962644579
$arr[17,2] = 5.888161900185501
$arr[17,3] = 11.69732181435834
$arr[17,4] = 27.30172781794434
$arr[17,5] = 0
$arr[17,6] = 7.344005520526261
$arr[17,7] = 24.43031386457932
$arr[17,8] = 0
$arr[17,9] = 9.71888413286889
$arr[17,10] = 9.820073821143753
$arr[17,11] = 0
$arr[17,12] = 0
$arr[17,13] = 24.40898293250917
$arr[18,0] = 12.70600227370191
$arr[18,1] = 8.296552048234298
$arr[18,2] = 5.90728040911685
$arr[18,3] = 11.69643134248035
$arr[18,4] = 27.29360750303287
$arr[18,5] = 0
$arr[18,6] = 7.344005520526261
$arr[18,7] = 24.39901937545069
$arr[18,8] = 0
$arr[18,9] = 9.800396757698584
$arr[18,10] = 9.829920881023337
$arr[18,11] = 0
$arr[18,12] = 0
$arr[18,13] = 24.38464303959133
$arr[19,0] = 13.08028946478006
$arr[19,1] = 8.335469547405197
$arr[19,2] = 5.970550673059481
$arr[19,3] = 11.69594758360909
$arr[19,4] = 27.27459620778698
$arr[19,5] = 0
$arr[19,6] = 7.344005520526261
$arr[19,7] = 24.29936317458361
$arr[19,8] = 0
$arr[19,9] = 10.07044666477051
$arr[19,10] = 9.864861536144479
$arr[19,11] = 0
$arr[19,12] = 0
$arr[19,13] = 24.30895021364378
$arr[20,0] = 13.32033338800919
$arr[20,1] = 8.360737183977918
$arr[20,2] = 6.011156412622159
$arr[20,3] = 11.69749570257589
$arr[20,4] = 27.26833984885148
$arr[20,5] = 0
$arr[20,6] = 7.344005520526261
$arr[20,7] = 24.23832966848188
$arr[20,8] = 0
$arr[20,9] = 10.24381234995751
$arr[20,10] = 9.889057568427765
$arr[20,11] = 0
$arr[20,12] = 0
$arr[20,13] = 24.26402754985255
$arr[21,0] = 13.19266038036997
$arr[21,1] = 8.347267433234904
$arr[21,2] = 5.989554158180232
$arr[21,3] = 11.69649921248997
$arr[21,4] = 27.27111444672213
$arr[21,5] = 0
$arr[21,6] = 7.344005520526261
$arr[21,7] = 24.27052904732896
$arr[21,8] = 0
$arr[21,9] = 10.1515872524861
$arr[21,10] = 9.876020738770501
$arr[21,11] = 0
$arr[21,12] = 0
$arr[21,13] = 24.28758687792504
$arr[22,0] = 12.69981315557991
$arr[22,1] = 8.295913309532116
$arr[22,2] = 5.906234077481242
$arr[22,3] = 11.69647079569441
$arr[22,4] = 27.29402232653077
$arr[22,5] = 0
$arr[22,6] = 7.344005520526261
$arr[22,7] = 24.40071728609043
$arr[22,8] = 0
$arr[22,9] = 9.795934109879864
$arr[22,10] = 9.829373089068952
$arr[22,11] = 0
$arr[22,12] = 0
$arr[22,13] = 24.38595680382783
$arr[23,0] = 12.15120890833922
$arr[23,1] = 8.239878446320169
$arr[23,2] = 5.813328195577405
$arr[23,3] = 11.70447453801876
$arr[23,4] = 27.34518447720579
$arr[23,5] = 0
$arr[23,6] = 7.344005520526261
$arr[23,7] = 24.55869670598306
$arr[23,8] = 0
$arr[23,9] = 9.400740579616746
$arr[23,10] = 9.785041076537905
$arr[23,11] = 0
$arr[23,12] = 0
$arr[23,13] = 24.51151349514885
$ws.Range("B2:O25").Value = $arr
Write-Output "done"
